# "correction \times & excel"
#
# Adds a LINEST() linear-regression array formula (slope/intercept plus the
# full regression statistics block) in H2:I6, based on G17:G38 (y) against
# C17:C38 (x); adds a helper ratio in I10 (=H2/H3); and repositions the
# existing chart so it no longer overlaps the new H:I columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LINEST full-statistics array formula spilling into H2:I6.
$ws.Range("H2:I6").FormulaArray = "=LINEST(G17:G38,C17:C38,TRUE,TRUE)"

# Ratio of slope to its standard error.
$ws.Range("I10").Formula = "=H2/H3"

# Move the chart two columns/rows to the right so it clears the new table.
$co = $ws.ChartObjects(1)
$co.Left = 581.1875
$co.Top = 43.75

# Match the author's final selection.
$ws.Range("I11").Select() | Out-Null
